$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (reuse existing header style/formatting from H1)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I column is always 1, J column mirrors H column values
$hValues = @(5, 6, 5, 3, 2)
for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValues[$i]
}
